$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained two more yearly columns (2021 and 2022), appended right
# after the existing last "2020" column (L). Mirror that column's
# formatting into the two new columns, for both the header row (3) and the
# data row (4).
$ws.Range("L3:L4").Copy()
$ws.Range("M3:M4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N3:N4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# New year headers.
$ws.Cells.Item(3, 13).Value = 2021
$ws.Cells.Item(3, 14).Value = 2022

# New data values, matching the existing figure (6.18) carried over from
# the 2019/2020 columns.
$ws.Cells.Item(4, 13).Value = 6.18
$ws.Cells.Item(4, 14).Value = 6.18

# The sheet's recorded selection moved on to N15 after the edit.
$ws.Range("N15").Select()
